$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "SGMTEST2402938482",
    "SGMTEST2208428758",
    "SGMTEST3559709487",
    "SGMTEST3938342818",
    "SGMTEST3585528276",
    "SGMTEST3132943337",
    "SGMTEST8815228500",
    "SGMTEST5936483766",
    "SGMTEST4621329996",
    "SGMTEST9085949196",
    "SGMTEST4069756425",
    "SGMTEST3850486410",
    "SGMTEST5761812024",
    "SGMTEST4047896363",
    "SGMTEST5142352881"
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

$ws.Columns.Item(2).ColumnWidth = 20.166666666666668

$ws.Range("D12").Select()
